# print_request_logsheet.xlsx edit
# - Fix incorrect material type (PLA -> Polylite) for July 2018 rows 20-26 (E column)
# - Fix incorrect material type (PLA -> Polylite) for August 2018 rows 2-7, 9-10 (E column)
# - Add a new request row (row 11) in August 2018 for the EVHP Holder TEE Track Concept 6 part
# - Update sheet selection / active tab state to match the edited state

$wb = $excel.ActiveWorkbook

$july = $wb.Worksheets.Item("July 2018")
$aug  = $wb.Worksheets.Item("August 2018")

# --- July 2018: correct material type PLA -> Polylite for rows 20-26 ---
foreach ($r in 20..26) {
    $july.Cells.Item($r, 5).Value = "Polylite"
}

# --- August 2018: correct material type PLA -> Polylite for rows 2-7 and 9-10 ---
foreach ($r in 2..7) {
    $aug.Cells.Item($r, 5).Value = "Polylite"
}
foreach ($r in 9..10) {
    $aug.Cells.Item($r, 5).Value = "Polylite"
}

# --- August 2018: add new row 11 for the new print request ---
# Column A holds a plain-text (not a real date serial) date string just like
# the other rows in this sheet, so force text entry with a leading apostrophe
# and then restore the normal (non quote-prefixed) cell formatting used by
# the rest of the column via a formats-only paste from an existing cell.
$aDest = $aug.Cells.Item(11, 1)
$aDest.Value = "'06-08-2018"
$aSrc = $aug.Cells.Item(9, 1)
$aSrc.Copy($null) | Out-Null
$aDest.PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$aug.Cells.Item(11, 3).Value = "EVHP Holder TEE Track Concept 6"
$aug.Cells.Item(11, 4).Value = 1
$aug.Cells.Item(11, 5).Value = "Polylite"
$aug.Cells.Item(11, 6).Value = 2
$aug.Cells.Item(11, 7).Value = 20
$aug.Cells.Item(11, 8).Value = 0.2

# --- View state: August 2018 selection moves to E7 ---
$aug.Activate()
$aug.Range("E7").Select() | Out-Null

# --- View state: July 2018 becomes the active tab, selection moves to E19:E26 ---
$july.Activate()
$july.Range("E19:E26").Select() | Out-Null
